$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range and swap columns C and D (codeforiati:group-code
# and codeforiati:group-name), including the header row, since the column
# order was flipped so that group-name now precedes group-code.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value()
    $dVal = $dCell.Value()
    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
